$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: "I have pushed my code to GitHub for grading (Check box if
# true)." + empty checkbox content control  ->  free-text "Write yes /
# if true): ____" answer (checkbox content control removed entirely).
# ---------------------------------------------------------------------

# Remove the checkbox content control (and its "☐" content) first.
for ($i = $d.ContentControls.Count; $i -ge 1; $i--) {
    $cc = $d.ContentControls.Item($i)
    if ($cc.Type -eq 8) {
        $cc.Delete($true)
    }
}

$oldPushed = "I have pushed my code to GitHub for grading (Check box if true).  "
$newPushed = "I have pushed my code to GitHub for grading (Write yes if true): ____"
$d.Content.Find.Execute($oldPushed, $true, $false, $false, $false, $false, $true, 1, $false, $newPushed, 2) | Out-Null

# ---------------------------------------------------------------------
# Edit 2: bump the four "Deliverable N" labels in the Measurement Data
# list by one (1->2, 2->3, 3->4, 4->5). Each is scoped with enough
# surrounding context to uniquely identify the right paragraph (there
# is an unrelated, earlier "Deliverable 1: Using KiCad..." paragraph
# that must stay untouched).
# ---------------------------------------------------------------------

$bumps = @(
    @{ old = "Deliverable 1: Quantized waveform of DAC";        new = "Deliverable 2: Quantized waveform of DAC" },
    @{ old = "Deliverable 2: Quantized waveform of ADC";        new = "Deliverable 3: Quantized waveform of ADC" },
    @{ old = "Deliverable 3: System module encoder";            new = "Deliverable 4: System module encoder" },
    @{ old = "Deliverable 4: Quantification of system";         new = "Deliverable 5: Quantification of system" }
)

foreach ($b in $bumps) {
    $d.Content.Find.Execute($b.old, $true, $false, $false, $false, $false, $true, 1, $false, $b.new, 2) | Out-Null
}
